# "added the pollen category viewer"
#
# The pollen data sample sheet gets trimmed down: a couple of sample data
# points that were only there for earlier testing are cleared out now that
# the category viewer reads the sheet directly, and the active selection is
# left on the cell the author was last working with (E4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# March / Earth sample value is no longer needed.
$ws.Range("C2").ClearContents()

# June row only keeps its "Garbage" category cell (now blank, formatting kept);
# the Grass/Oil sample values for that row are removed.
$ws.Range("B4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

# Leave the selection where the author last left it.
$ws.Range("E4").Select()
